# "screening finished by me" - Jerome (Reviewer 1) fills in his Decision/Notes
# columns (B/C) on the "Reviewer1" sheet for the papers he has screened
# (rows 3-48), mirroring the same per-paper decisions that already exist on
# the "Reviewer2" sheet's columns A/B for those rows.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Reviewer1")
$ws2 = $wb.Worksheets.Item("Reviewer2")

# row -> Decision (column B), Notes (column C, optional)
$decisions = @{
    3  = @("ok", "predictive model")
    4  = @("ok", $null)
    5  = @("ok", "simulating BP models")
    6  = @("ok", $null)
    7  = @("ok", "von uns :D")
    8  = @("ok", "Daniels Diss")
    9  = @("ok", $null)
    10 = @("no DT", $null)
    11 = @("ok", $null)
    12 = @("ok", $null)
    13 = @("ok", "nicht sicher ob wirklich MDE")
    14 = @("secondary study", $null)
    15 = @("ok", "no DT?")
    16 = @("ok", "no DT?")
    17 = @("no MDE", $null)
    18 = @("ok", $null)
    19 = @("ok", $null)
    20 = @("ok", $null)
    21 = @("no DT", $null)
    22 = @("no DT", $null)
    23 = @("ok", "sie haben ein Modell und nennen ihr modellverarbeitendes Tool DT ")
    24 = @("ok", $null)
    25 = @("experience", $null)
    26 = @("ok", $null)
    27 = @("ok", $null)
    28 = @("ok", $null)
    29 = @("no MDE", $null)
    30 = @("no MDE", "mathematical model")
    31 = @("experience", $null)
    32 = @("no MDE", $null)
    33 = @("challenges", $null)
    34 = @("Titel vs. Abstract?", $null)
    35 = @("no MDE, no DT", $null)
    36 = @("no DT", $null)
    37 = @("ok", $null)
    38 = @("no DT", $null)
    39 = @("ok", $null)
    40 = @("ok", "no DT?")
    41 = @("ok", $null)
    42 = @("ok", $null)
    43 = @("no DT", $null)
    44 = @("ok", "no DT?")
    45 = @("no MDE + no DT", $null)
    46 = @("ok", "no DT?")
    47 = @("experience", $null)
    48 = @("ok", "no DT?")
}

for ($row = 3; $row -le 48; $row++) {
    $pair = $decisions[$row]
    $ws1.Cells.Item($row, 2).Value2 = $pair[0]
    if ($pair[1]) {
        $ws1.Cells.Item($row, 3).Value2 = $pair[1]
    }
}

# Window/view bookkeeping to match where the reviewer ended up scrolled to
# when they saved: Reviewer2 had been scrolled near the bottom (topLeftCell
# A89 / C102 selected) and is now back near the top; Reviewer1 is zoomed in
# and scrolled down to the rows just screened.
$ws2.Activate()
$ws2.Range("A8").Select()

$ws1.Activate()
$ws1.Range("C49").Select()
$excel.ActiveWindow.Zoom = 130
